$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Controller")
$ws.Range("E3").Copy() | Out-Null
$ws.Range("F3").PasteSpecial(-4122) | Out-Null
$ws.Range("F3").Value = $ws.Range("E3").Value2
